$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.619.99"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.444.45"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'573.07"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'158.94"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.443.88"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  -6.62%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "4.036.29"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'27.48"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("E16").Value = "  -9.25%  "
$ws.Range("D17").Value = "64.641.46"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "3.447.84"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").Value = "'13.76"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "'380.23"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'7.96"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'72.35"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'0.530"
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("D26").Value = "'0.0000118"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'1.02"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.09"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "'2.00"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'23.22"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'161.38"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").Value = "2.880.27"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'0.0747"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "'26.28"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").Value = "'0.797"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'42.95"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'25.91"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0310"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'2.44"
$ws.Range("E47").Value = "  +13.10%  "
$ws.Range("D48").Value = "'322.79"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'0.844"
$ws.Range("E51").Value = "  -2.51%  "
